$d = $word.ActiveDocument

# --- Simple single-run text replacements (formatting preserved automatically) ---

$r = $d.Content
$r.Find.Execute("Unveiling the Quantum Realm: A Path to Technological Singularity") | Out-Null
$r.Text = "Mathematics: A Journey Through Numbers and Logic"

$r = $d.Content
$r.Find.Execute("Dr") | Out-Null
$r.Text = "Prof"

$r = $d.Content
$r.Find.Execute(" Amelia Jacobson") | Out-Null
$r.Text = " Marcus Adams"

$r = $d.Content
$r.Find.Execute("amelia") | Out-Null
$r.Text = "marcus"

$r = $d.Content
$r.Find.Execute("jacobson@quantuminstitute") | Out-Null
$r.Text = "adams@eduinstitute"

# --- Paragraph 5 (main body paragraph): several sentence rewrites, plus two
#     sentences split into extra runs (new "." run + new sentence run) ---

$p5 = $d.Paragraphs(5).Range
$p5.InsertXML('<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00EA4683" w:rsidRDefault="00C50137"><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>Mathematics, the language of the universe, holds immense power to unveil the secrets of the cosmos</w:t></w:r><w:r w:rsidR="0022190B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> It is a subject that captivates the minds of young and old alike, enriching our understanding of nature and enabling us to unravel the mysteries it holds</w:t></w:r><w:r w:rsidR="0022190B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> From the intricate patterns of fractals to the elegance of calculus, mathematics invites us on an exploration of the fundamental principles governing our world</w:t></w:r><w:r w:rsidR="0022190B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>In the realm of mathematics, we investigate the complex interplay of numbers, symbols, and equations</w:t></w:r><w:r w:rsidR="0022190B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> We unlock the mysteries of geometry, traversing through the landscapes of shapes and angles</w:t></w:r><w:r w:rsidR="0022190B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> We discover the beauty of algebra, manipulating expressions and equations to unveil hidden relationships</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> And as we delve into calculus, we uncover the intricate dance of change, exploring the rates and patterns that define the universe''s dynamic processes</w:t></w:r><w:r w:rsidR="0022190B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>Mathematics is not merely about abstract concepts; it is a tool that empowers us to solve real-world problems</w:t></w:r><w:r w:rsidR="0022190B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> From predicting the trajectory of a rocket to analyzing financial data, mathematics equips us with the skills to navigate the complexities of our world</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> It fosters critical thinking, analytical reasoning, and problem-solving abilities, preparing us for success in various fields and endeavors</w:t></w:r><w:r w:rsidR="0022190B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# --- Paragraph 7 (summary body paragraph): sentence rewrites, removal of the
#     lastRenderedPageBreak run (merged away), two new split runs, and a new
#     trailing empty paragraph appended at the very end of the document ---

$p7 = $d.Paragraphs(7).Range
$p7.InsertXML('<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00EA4683" w:rsidRDefault="00C50137"><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>Mathematics stands as a powerful and versatile tool that unveils the mysteries of the universe and equips us with essential skills for navigating the complexities of life</w:t></w:r><w:r w:rsidR="0022190B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> It invites us on an exploration of numbers, symbols, and equations, captivating our minds with its intricate patterns and elegant structures</w:t></w:r><w:r w:rsidR="0022190B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> As we delve deeper into mathematics, we unlock the secrets of geometry, algebra, and calculus, gaining a profound understanding of the fundamental principles that govern our world</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> Mathematics is more than a subject; it is a language that empowers us to solve real-world problems and make informed decisions, shaping our lives and contributing to the advancement of society</w:t></w:r><w:r w:rsidR="0022190B"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
